$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: item_count changes from 1 to 10
$ws.Range("E3").Value = 10

# New row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "event"
$ws.Range("C4").Value = 14
$ws.Range("D4").Value = 100002
$ws.Range("E4").Value = 10000

# New row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "event"
$ws.Range("C5").Value = 14
$ws.Range("D5").Value = 100003
$ws.Range("E5").Value = 10

# Update selection to E5, matching the saved view state
[void]$ws.Range("E5").Select()
